# Edit the "backlog" workbook:
#  - Sheet "Hoja1": fix a typo in C3, add a new description in C4, and
#    append three new backlog rows (15, 16, 17) with a new "leader board" /
#    scoring feature.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Fix typo: "preguntar" -> "preguntas"
$ws.Range("C3").Value = "lógica de siguiente pregunta, guardar en json, que haya preguntas de distinta dificultad"

# New description text for "Entender respuestas"
$ws.Range("C4").Value = "A partir de los response de los jugadores asignar los puntajes correspondientes"

# New backlog rows
$ws.Range("A15").Value = 14
$ws.Range("B15").Value = "Leader board"
$ws.Range("C15").Value = "Muestra los puntajes parciales de todos los jugadores"

$ws.Range("A16").Value = 15
$ws.Range("B16").Value = "Agregar los puntajes a los jugadores"
$ws.Range("C16").Value = "Cuando se agrega un jugador empieza con 0. Y se pueden sumar o restar"

$ws.Range("A17").Value = 16
$ws.Range("B17").Value = "Enviarle preguntas a los jugadores"
$ws.Range("C17").Value = "Enviarles un request con la pregunta a los jugadores"

# Column widths widened to best-fit the new, longer content
$ws.Columns.Item(2).ColumnWidth = 32.47
$ws.Columns.Item(3).ColumnWidth = 83.35

# Move the active selection to B4
$ws.Range("B4").Select() | Out-Null
